$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.855.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.47%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.089.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'526.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.24%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -3.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.088.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.53%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.06%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.618.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.91%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.88%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'25.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -7.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -2.51%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'57.846.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.45%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.077.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.71%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.60%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -4.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'342.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.86%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.55%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'67.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.85%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -1.83%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.53%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.38%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'20.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.10%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'158.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.91%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.50%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'25.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.53%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.51%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0666"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.52%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.94%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.20%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.683"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.127.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'36.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.0262"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.277.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.39%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.72%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'20.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.96%  "
$ws.Range("E51").Style = "Normal"
